$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update version number from 0.1 to 1.0 (must remain text, not be coerced to a number)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1.0"

# Swap the TC3 / TC4 step & result content.
$atribuirStep   = "Chefe Dado um registro selecionado (solicitação aguardando autorização de pagamento - AP), o usuário pode atribuir/desatribuir a responsabilidade da AP a si próprio; e Clica para atribuir/desatribuir o registro a si mesmo."
$atribuirResult = "SYSTEM Atualiza a lista de registros de solicitações, onde o nome deverá constar o nome do usuário logado (que se atribuiu como responsável pela AP) no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."
$filtrarStep    = "Chefe Seleciona um usuário para filtrar as autorizações de pagamento associadas a ele; e Submete a busca ao sistema."
$filtrarResult  = "SYSTEM Filtra os registros (autorizações de pagamento pendentes) e exibe apenas aqueles atribuídos ao usuário selecionado."

# Row 28 (TC3's 2nd step row) previously held the "atribuir" content; it should now hold the "filtrar" content.
$ws.Range("B28").Value = $filtrarStep
$ws.Range("D28").Value = $filtrarResult

# Row 36 (TC4's 2nd step row) previously held the "filtrar" content; it should now hold the "atribuir" content.
$ws.Range("B36").Value = $atribuirStep
$ws.Range("D36").Value = $atribuirResult
